# Update LipidLynxX_test.csv for GitHubActions
# Shift the SOURCE_05 / SOURCE_05_converted values (columns I and J) for
# rows 8-17 down by one row, and place a new truncated value in row 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing values first (Value2 reads correctly; Value getter
# is unreliable on this runtime), then write them back shifted down by one
# row so that earlier writes don't clobber values still to be read.
$i8 = $ws.Range("I8").Value2
$i9 = $ws.Range("I9").Value2
$i10 = $ws.Range("I10").Value2
$i11 = $ws.Range("I11").Value2
$i12 = $ws.Range("I12").Value2
$i13 = $ws.Range("I13").Value2
$i14 = $ws.Range("I14").Value2
$i15 = $ws.Range("I15").Value2
$i16 = $ws.Range("I16").Value2

$j8 = $ws.Range("J8").Value2
$j9 = $ws.Range("J9").Value2
$j10 = $ws.Range("J10").Value2
$j11 = $ws.Range("J11").Value2
$j12 = $ws.Range("J12").Value2
$j13 = $ws.Range("J13").Value2
$j14 = $ws.Range("J14").Value2
$j15 = $ws.Range("J15").Value2
$j16 = $ws.Range("J16").Value2

$ws.Range("I17").Value = $i16
$ws.Range("J17").Value = $j16

$ws.Range("I16").Value = $i15
$ws.Range("J16").Value = $j15

$ws.Range("I15").Value = $i14
$ws.Range("J15").Value = $j14

$ws.Range("I14").Value = $i13
$ws.Range("J14").Value = $j13

$ws.Range("I13").Value = $i12
$ws.Range("J13").Value = $j12

$ws.Range("I12").Value = $i11
$ws.Range("J12").Value = $j11

$ws.Range("I11").Value = $i10
$ws.Range("J11").Value = $j10

$ws.Range("I10").Value = $i9
$ws.Range("J10").Value = $j9

$ws.Range("I9").Value = $i8
$ws.Range("J9").Value = $j8

$ws.Range("I8").Value = "Cer(d18:0/26:0)"
$ws.Range("J8").Value = "Cer(18:0/26:0)"
